$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 6, shifting the existing row 6
# ("even_MAG-GUT74962.fa") down to row 7.
$ws.Rows("6:6").Insert()

# Match the style of the other data rows in column A (border/alignment)
# by copying formatting (not value) from the row above.
$ws.Range("A5").Copy()
$ws.Range("A6").PasteSpecial(-4122)

# Populate the newly inserted row 6 with the new record.
$ws.Range("A6").Value = "even_MAG-GUT72020.fa"
$ws.Range("B6").Value = -0.8909409175172307
$ws.Range("C6").Value = "s__CAG-194 sp000432915"
$ws.Range("D6").Value = "s__CAG-194 sp000432915(reject)"
